$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B4").Value = 155680.43137199999
$ws.Range("A42").Value = 155680.43137199999

$ws.Range("B41").Value = 1.5747536799999999
$ws.Range("C41").Value = 2.0575342000000001
$ws.Range("D41").Value = 2.6258807000000002
$ws.Range("E41").Value = 3.14089441
$ws.Range("F41").Value = 3.6407742399999998
$ws.Range("G41").Value = 5.1078768099999996
$ws.Range("H41").Value = 10.57204033

$ws.Range("H46").Select()
